# Remove the "Soustředění na hru" entries: the Czech diagnostics sheet
# had three helper rows (398:400) that all referenced that label; the
# author deleted those whole rows, which shifts every row below them
# up by three and drops the now-unused shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("398:400").Delete() | Out-Null

# Restore the cursor/viewport roughly where the author left it after
# performing the deletion.
$ws.Range("B402").Select()
